# Update "想去人数" (want-to-go count) figures on both the "展览" sheet and
# the aggregated "全部类型" sheet to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 542
$ws1.Range("F12").Value = 13378
$ws1.Range("F13").Value = 169
$ws1.Range("F16").Value = 5501

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F31").Value = 542
$ws4.Range("F34").Value = 13378
$ws4.Range("F35").Value = 169
$ws4.Range("F39").Value = 5501
